# Reorders the "Periodo Mora" / "Valor Mora" table (rows 16-50) from newest-first
# to oldest-first chronological order, and carries each period's "Valor Mora"
# amount along with it (period 2009 keeps its distinctive 52267 value, which
# now lands on row 50 instead of row 16; every other period uses 56000).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @(
    "1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003","2004","2005","2006","2007","2008","2009"
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("E$row").Value = $periods[$i]
    if ($periods[$i] -eq "2009") {
        $ws.Range("F$row").Value = 52267
    } else {
        $ws.Range("F$row").Value = 56000
    }
}
